$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 372, shifting existing rows 372:490 down to 373:491
$ws.Rows.Item(372).Insert()

# Populate the newly inserted row 372 with the new weekly record
$ws.Range("A372").Value = 4
$ws.Range("B372").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C372").Value = "Los Lagos"
$ws.Range("D372").Value = 44985
$ws.Range("E372").Value = 10
$ws.Range("F372").Value = 100112008
$ws.Range("G372").Value = "Coliflor"
$ws.Range("H372").Value = "Sin especificar"
$ws.Range("I372").Value = "Primera"
$ws.Range("J372").Value = 1000
$ws.Range("K372").Value = 1700
$ws.Range("L372").Value = 1700
$ws.Range("M372").Value = 1700
$ws.Range("N372").Value = "$/unidad"
$ws.Range("O372").Value = "Región Metropolitana"
$ws.Range("P372").Value = 1700
$ws.Range("Q372").Value = 1
$ws.Range("R372").Value = "Hortaliza"
